# Daily attendance processing - 2026-01-03 09:34:06
# Swap the "Recorded By" text from "System, dnasr281@gmail.com" to
# "dnasr281@gmail.com, System" for the specific session rows touched by
# today's processing run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,7,16,17,21,22,36,37,41,42,56,57,61,62,76,77,81,82,83,84,85,86,95,96,100,101,102,103,104,105,114,115,119,120,121,122,123,124,133,134,138,139,140,141,142,143,152,153,157,158,159,160,161,162,171,172,176,177,191,192,196,197,211,212,216,217,231,232)

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

foreach ($r in $rows) {
    $cell = $ws.Range("G$r")
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}

Write-Host "Updated $($rows.Count) cells in column G"
